$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-27 Monday", "2025-01-28 Tuesday"),
    @("895÷2=", "231÷6="),
    @("214÷4=", "331÷3="),
    @("158÷4=", "872÷3="),
    @("429÷2=", "690÷9="),
    @("736÷6=", "791÷2="),
    @("215÷2=", "934÷2="),
    @("411÷9=", "718÷4="),
    @("834÷8=", "817÷7="),
    @("999÷6=", "334÷6="),
    @("515÷3=", "732÷6="),
    @("269÷4=", "475÷2="),
    @("564÷8=", "458÷4="),
    @("103÷6=", "433÷3="),
    @("682÷8=", "859÷4="),
    @("149÷8=", "370÷9="),
    @("879÷6=", "683÷6="),
    @("253÷8=", "801÷5="),
    @("106÷8=", "439÷5="),
    @("511÷9=", "941÷2="),
    @("554÷5=", "742÷4="),
    @("583÷9=", "480÷5="),
    @("614÷2=", "145÷5="),
    @("606÷8=", "883÷9="),
    @("170÷9=", "577÷8="),
    @("919÷5=", "425÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
